# ==========================================================================
# [ADDITIONAL SCRAPING] added code to scrape more data about a player's
# batting performance in a match, also updated the excel sheets
#
# This script:
#   1. Adds a new "Player Info" sheet (first tab) with player bio data.
#   2. Renames MATCH_CARD_LINK -> MATCH_CODE on "ODI Batting" and replaces
#      the full howstat.com scorecard URLs with the bare numeric match
#      code in column D.
#   3. Does the same MATCH_CARD_LINK -> MATCH_CODE rename + URL -> code
#      replacement on "ODI Bowling" (column B there).
#   4. Adds a new "ODI Batting Extra" sheet (last tab) with additional
#      per-match batting stats.
# ==========================================================================

$wb = $excel.ActiveWorkbook

function Set-HeaderCell($cell) {
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108   # xlCenter
    $cell.VerticalAlignment = -4160     # xlTop
    $cell.Borders.LineStyle = 1
}

# Writing a leading-apostrophe value (eg "'4186") is how Excel forces a
# numeric-looking literal to be stored as text instead of a number. Excel
# marks such cells with a "quote prefix" style internally; resetting the
# style back to Normal afterwards keeps the cell's type as Text while
# dropping that incidental formatting so the cell matches a plain
# text cell with no special style.
function Set-TextValue($cell, $value) {
    $cell.Value = "'" + $value
    $cell.Style = "Normal"
}

# --------------------------------------------------------------------------
# 1) New "Player Info" sheet, inserted before the current first sheet so it
#    becomes the new first tab.
# --------------------------------------------------------------------------
$firstSheet = $wb.Worksheets.Item(1)
$playerInfo = $wb.Worksheets.Add($firstSheet)
$playerInfo.Name = "Player Info"

$playerInfoHeaders = @("ID", "NAME", "BATTING_HAND", "BOWL_STYLE")
for ($c = 1; $c -le $playerInfoHeaders.Length; $c++) {
    $cell = $playerInfo.Cells.Item(1, $c)
    $cell.Value = $playerInfoHeaders[$c - 1]
    Set-HeaderCell $cell
}

Set-TextValue $playerInfo.Cells.Item(2, 1) "4340"
$playerInfo.Cells.Item(2, 2).Value = "Reeza Raphael Hendricks"
$playerInfo.Cells.Item(2, 3).Value = "Right Handed"
$playerInfo.Cells.Item(2, 4).Value = "Right Arm Medium Fast"

# --------------------------------------------------------------------------
# 2) "ODI Batting": header rename + URL -> bare match code in column D.
# --------------------------------------------------------------------------
$odiBatting = $wb.Worksheets.Item("ODI Batting")
$odiBatting.Range("D1").Value = "MATCH_CODE"

$battingCodes = @("4186","4187","4188","4206","4207","4208","4222","4224","4226","4237","4238","4241","4244","4247","4261","4264","4269","4271","4401","4405","4408","4488","4491","4517","4657","4658","4700")
for ($i = 0; $i -lt $battingCodes.Length; $i++) {
    $row = $i + 2
    Set-TextValue $odiBatting.Cells.Item($row, 4) $battingCodes[$i]
}

# --------------------------------------------------------------------------
# 3) "ODI Bowling": header rename + URL -> bare match code in column B.
# --------------------------------------------------------------------------
$odiBowling = $wb.Worksheets.Item("ODI Bowling")
$odiBowling.Range("B1").Value = "MATCH_CODE"

$bowlingCodes = @("4237","4241","4244")
for ($i = 0; $i -lt $bowlingCodes.Length; $i++) {
    $row = $i + 2
    Set-TextValue $odiBowling.Cells.Item($row, 2) $bowlingCodes[$i]
}

# --------------------------------------------------------------------------
# 4) New "ODI Batting Extra" sheet, appended after the current last sheet.
# --------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$battingExtra = $wb.Worksheets.Add($null, $lastSheet)
$battingExtra.Name = "ODI Batting Extra"

$extraHeaders = @("MATCH_CODE", "BATTING_POSITION", "NUM_4", "NUM_6", "PERCENT_RUNS_OF_TOTAL", "MAN_OF_MATCH")
for ($c = 1; $c -le $extraHeaders.Length; $c++) {
    $cell = $battingExtra.Cells.Item(1, $c)
    $cell.Value = $extraHeaders[$c - 1]
    Set-HeaderCell $cell
}

# MATCH_CODE, BATTING_POSITION, NUM_4, NUM_6, PERCENT_RUNS_OF_TOTAL, MAN_OF_MATCH
# "EMPTY" is a sentinel meaning "leave this cell blank" (rows with no data).
$extraData = @(
    @("4224", "2", "1", "0", "7.14%",  "NO"),
    @("4226", "2", "2", "0", "2.50%",  "NO"),
    @("4237", "2", "5", "0", "16.92%", "NO"),
    @("4238", "1", "1", "0", "2.42%",  "NO"),
    @("4241", "3", "8", "2", "44.39%", "YES"),
    @("4244", "3", "0", "0", "1.22%",  "NO"),
    @("4247", "EMPTY", "EMPTY", "EMPTY", "EMPTY", "NO"),
    @("4261", "2", "0", "0", "0.43%",  "NO"),
    @("4264", "2", "4", "0", "11.55%", "NO"),
    @("4269", "2", "0", "0", "1.21%",  "NO"),
    @("4271", "EMPTY", "EMPTY", "EMPTY", "EMPTY", "NO"),
    @("4401", "EMPTY", "EMPTY", "EMPTY", "EMPTY", "NO"),
    @("4405", "EMPTY", "EMPTY", "EMPTY", "EMPTY", "NO"),
    @("4408", "EMPTY", "EMPTY", "EMPTY", "EMPTY", "NO"),
    @("4488", "3", "5", "0", "18.02%", "NO"),
    @("4491", "3", "0", "0", "0.80%",  "NO"),
    @("4517", "EMPTY", "EMPTY", "EMPTY", "EMPTY", "NO"),
    @("4657", "3", "9", "1", "26.62%", "NO"),
    @("4658", "3", "0", "0", "3.03%",  "NO"),
    @("4700", "2", "6", "0", "18.12%", "NO")
)

for ($i = 0; $i -lt $extraData.Length; $i++) {
    $row = $i + 2
    $rowdata = $extraData[$i]

    Set-TextValue $battingExtra.Cells.Item($row, 1) $rowdata[0]

    if ($rowdata[1] -eq "EMPTY") {
        # leave BATTING_POSITION / NUM_4 / NUM_6 / PERCENT_RUNS_OF_TOTAL blank
    } else {
        $battingExtra.Cells.Item($row, 2).Value = [int]$rowdata[1]
        Set-TextValue $battingExtra.Cells.Item($row, 3) $rowdata[2]
        Set-TextValue $battingExtra.Cells.Item($row, 4) $rowdata[3]
        Set-TextValue $battingExtra.Cells.Item($row, 5) $rowdata[4]
    }

    $battingExtra.Cells.Item($row, 6).Value = $rowdata[5]
}
